$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Not worse",
    "Not worse",
    "Not worse",
    "Not worse",
    "A little worse",
    "Not worse",
    "Not worse",
    "Not worse",
    "A little worse",
    "Not worse",
    "Not worse",
    "Not worse",
    "Not worse",
    "Not worse"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

$ws.Range("D16").Select()
